$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.868.45"
$ws.Range("E2").Value = "  -0.99%  "

$ws.Range("D3").Value = "2.512.05"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.21%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.517"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").Value = "2.512.81"
$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.83%  "

$ws.Range("E11").Value = "  -0.94%  "

$ws.Range("E12").Value = "  +3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "

$ws.Range("D14").Value = "2.983.58"
$ws.Range("E14").Value = "  -1.01%  "

$ws.Range("D15").Value = "69.761.60"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.92"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "2.537.54"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "349.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.27%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.88%  "

$ws.Range("D28").Value = "2.646.93"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "0.0₃0899"
$ws.Range("E30").Value = "  -2.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "464.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.12%  "

$ws.Range("E34").Value = "  -1.52%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.60%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("E43").Value = "  -3.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.29%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0732"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.40%  "
